$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old "mortgage" (row 5) and "recession" (row 6) rows.
# Delete row 6 first so row indices for the remaining deletion stay correct.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Insert a new row above row 2 to make room for "bank"
$ws.Rows.Item(2).Insert()

# Row 2: bank
$ws.Cells.Item(2, 1).Value = "bank"
$ws.Cells.Item(2, 2).Value = "Frequency"
$ws.Cells.Item(2, 3).Value = 0.5336

# Match formatting of the other data rows (border + bold + centered for A/B,
# plain/default for C) since the freshly inserted row otherwise gets a
# slightly different auto-generated style.
$ws.Cells.Item(3, 1).Copy() | Out-Null
$ws.Cells.Item(2, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 2).Copy() | Out-Null
$ws.Cells.Item(2, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(2, 3).ClearFormats() | Out-Null
